# Demo - Comments only
#
# Slide 3 ("Design") repositions a handful of shapes:
#   - the Title placeholder gets an explicit position/size (previously it
#     just inherited the layout's xfrm because <p:spPr/> was empty)
#   - "Flowchart: Magnetic Disk 4" (AZURE Cloud) moves up/left
#   - "Flowchart: Process 7" (SAP Solution Manager) moves up/left
#
# The PowerPoint object model stores Left/Top/Width/Height in points while
# the OOXML stores EMU (1 pt = 12700 EMU). A tiny (+0.5 EMU) nudge is added
# before converting to points so the internal float round-trip lands on the
# exact target EMU value instead of one EMU short.

function EMUToPt($emu) {
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Title placeholder: give it an explicit xfrm -----------------------
$title = $s.Shapes.Item("Title 1")
$title.Left   = EMUToPt 100705
$title.Top    = EMUToPt 227082
$title.Width  = EMUToPt 8765651
$title.Height = EMUToPt 310500

# --- "Flowchart: Magnetic Disk 4" (AZURE Cloud) -------------------------
$magDisk = $s.Shapes.Item("Flowchart: Magnetic Disk 4")
$magDisk.Left = EMUToPt 2334890
$magDisk.Top  = EMUToPt 3424174

# --- "Flowchart: Process 7" (SAP Solution Manager) ----------------------
$sapProcess = $s.Shapes.Item("Flowchart: Process 7")
$sapProcess.Left = EMUToPt 4363906
$sapProcess.Top  = EMUToPt 781834
